$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1. Insert a new data row (row 33) below the current last data row
#    (row 32), preserving the "closing" border style that belonged
#    to the old last row, then restore the "normal" style on row 32
#    by copying the formatting from row 31.
# ------------------------------------------------------------------
$ws.Range("B33:J33").Insert(-4121)                         # xlShiftDown
$ws.Range("B32:J32").Copy($ws.Range("B33:J33"))             # row33 <- old row32 (style + data, temp)
$ws.Range("B31:J31").Copy($ws.Range("B32:J32"))             # row32 <- row31 (style + data, temp)
$excel.CutCopyMode = 0

# ------------------------------------------------------------------
# 2. Write the final values for the worker detail table (rows 16-33)
#    Periods are now listed in ascending (chronological) order.
# ------------------------------------------------------------------

# -- Karen Margarita Mendoza Hernandez (rows 16-17) --
$ws.Cells.Item(16,2).Value = "CC"
$ws.Cells.Item(16,3).Value = "1047400498"
$ws.Cells.Item(16,4).Value = "KAREN MARGARITA MENDOZA HERNANDEZ"
$ws.Cells.Item(16,5).Value = "1806"
$ws.Cells.Item(16,6).Value = 26041
$ws.Cells.Item(16,7).Value = 900000

$ws.Cells.Item(17,2).Value = "CC"
$ws.Cells.Item(17,3).Value = "1047400498"
$ws.Cells.Item(17,4).Value = "KAREN MARGARITA MENDOZA HERNANDEZ"
$ws.Cells.Item(17,5).Value = "1807"
$ws.Cells.Item(17,6).Value = 31249
$ws.Cells.Item(17,7).Value = 900000

# -- Javier Antonio Blanco De La Rosa (rows 18-33) --
$periods = @("2405","2406","2407","2408","2409","2410","2411","2412","2501","2502","2503","2504","2505","2506","2507","2508")
$values  = @(41600,52000,52000,52000,52000,52000,52000,52000,52000,52000,52000,52000,52000,52000,52000,52000)

for ($i = 0; $i -lt $periods.Length; $i++) {
    $r = 18 + $i
    $ws.Cells.Item($r,2).Value = "CC"
    $ws.Cells.Item($r,3).Value = "1007974773"
    $ws.Cells.Item($r,4).Value = "JAVIER ANTONIO BLANCO DE LA ROSA"
    $ws.Cells.Item($r,5).Value = $periods[$i]
    $ws.Cells.Item($r,6).Value = $values[$i]
    $ws.Cells.Item($r,7).Value = 1300000
}

# ------------------------------------------------------------------
# 3. Update the summary fields: total "Valor Mora" and "Cant. Periodos"
# ------------------------------------------------------------------
$ws.Range("E11").Value = 878890
$ws.Range("F13").Value = 18

Write-Host "Done."
